$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.00002074986032285508, 0.00007097389502863649, 26.21740644021617, 8.660232485948974, 0, 34.87773064992049)
    3 = @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.649481363816475)
    4 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 12.59312877619104)
    5 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    6 = @(0.01514828764759746, 114.8270160096505, 26.21740644021617, 645.3272768299601, 1, 786.3868475674743)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    8 = @(0.6753301551942219, 10.29869402782916, 26.21740644021617, 8.660232485948974, 1, 45.85166310918853)
    9 = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 0, 4.837881874639075)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

$wb.Save()
